$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8
$ws.Range("A8").Value = -1450.392252482472
$ws.Range("B8").Value = -1450.392252482472
$ws.Range("C8").Value = -1450.392252482472
$ws.Range("D8").Value = -1339.250701011006
$ws.Range("E8").Value = -1339.250701011006
$ws.Range("F8").Value = -1034.248189140555
$ws.Range("G8").Value = -1247.897379399485
$ws.Range("H8").Value = -1225.78532072732
$ws.Range("I8").Value = -971.5862553205345
$ws.Range("J8").Value = -845.5579614850009
$ws.Range("K8").Value = -845.5579614850009
$ws.Range("L8").Value = -822.3941131771093
$ws.Range("M8").Value = -822.3941131771093
$ws.Range("N8").Value = -809.1449917242714
$ws.Range("O8").Value = -809.1449917242714
$ws.Range("P8").Value = -809.1449917242714
$ws.Range("Q8").Value = -774.5200934940744
$ws.Range("R8").Value = -774.5200934940744
$ws.Range("S8").Value = -774.5200934940744
$ws.Range("T8").Value = -774.5200934940744
$ws.Range("U8").Value = -606.1629503574208
$ws.Range("V8").Value = -606.1629503574208
$ws.Range("W8").Value = -606.1629503574208
$ws.Range("X8").Value = -606.1629503574208
$ws.Range("Y8").Value = -606.1629503574208
$ws.Range("Z8").Value = -606.1629503574208
$ws.Range("AA8").Value = -606
$ws.Range("AB8").Value = 'passed'
$ws.Range("AC8").Value = 'sub(protectedDiv(add(protectedDiv(add(x3, y3), y2), conditional(y2, y1)), y3), tan(ang_vel(x3, y1, y3, x1)))'
$ws.Range("AD8").Value = 'max, atan, limit, acos, x2, cos, asin, sin, '
$ws.Range("AE8").Value = 'tan, y2, ang_vel, sub, protectedDiv, conditional, y1, x1, y3, x3, add, '
$ws.Range("AF8").Value = 'replace'

# Row 9
$ws.Range("A9").Value = -1846.511939803875
$ws.Range("B9").Value = -1814.824988986811
$ws.Range("C9").Value = -1643.775424731376
$ws.Range("D9").Value = -1578.434304670724
$ws.Range("E9").Value = -1611.060645760238
$ws.Range("F9").Value = -1571.749795489733
$ws.Range("G9").Value = -1562.356788893979
$ws.Range("H9").Value = -1514.481553625497
$ws.Range("I9").Value = -1471.719288795401
$ws.Range("J9").Value = -839.9296345291842
$ws.Range("K9").Value = -555.7937755520754
$ws.Range("L9").Value = -555.7937755520754
$ws.Range("M9").Value = -555.7937755520754
$ws.Range("N9").Value = -491.9180312454844
$ws.Range("O9").Value = -491.9180312454844
$ws.Range("P9").Value = -368.582643133616
$ws.Range("Q9").Value = -491.9180312454844
$ws.Range("R9").Value = -470.7954276157421
$ws.Range("S9").Value = -422.3190395215176
$ws.Range("T9").Value = -390.5510587169434
$ws.Range("U9").Value = -390.5510587169434
$ws.Range("V9").Value = -390.5510587169434
$ws.Range("W9").Value = -386.6780925103196
$ws.Range("X9").Value = -339.7634857551262
$ws.Range("Y9").Value = -339.7634857551262
$ws.Range("Z9").Value = -339.7634857551262
$ws.Range("AA9").Value = -339
$ws.Range("AB9").Value = 'passed'
$ws.Range("AC9").Value = 'ang_vel(add(atan(protectedDiv(limit(x1, x1, x1), cos(y3)), y1), acos(y3, x2)), tan(ang_vel(x1, x1, y1, x2)), tan(atan(x2, x3)), sub(y3, atan(limit(x1, y2, x3), y1)))'
$ws.Range("AD9").Value = 'max, conditional, asin, sin, '
$ws.Range("AE9").Value = 'tan, ang_vel, y2, atan, protectedDiv, limit, y1, acos, x2, sub, x1, y3, x3, add, cos, '
$ws.Range("AF9").Value = 'replace'

# Row 10
$ws.Range("A10").Value = -1789.86075472485
$ws.Range("B10").Value = -1528.037694102767
$ws.Range("C10").Value = -1528.037694102767
$ws.Range("D10").Value = -1357.045211674417
$ws.Range("E10").Value = -1357.045211674417
$ws.Range("F10").Value = -1161.658315096036
$ws.Range("G10").Value = -1161.658315096036
$ws.Range("H10").Value = -1161.658315096036
$ws.Range("I10").Value = -1153.448653591956
$ws.Range("J10").Value = -1136.089205368618
$ws.Range("K10").Value = -1136.089205368618
$ws.Range("L10").Value = -1136.089205368618
$ws.Range("M10").Value = -1120.588294685263
$ws.Range("N10").Value = -1120.588294685263
$ws.Range("O10").Value = -1097.867324816126
$ws.Range("P10").Value = -1071.788871636817
$ws.Range("Q10").Value = -1066.342320838975
$ws.Range("R10").Value = -1040.95524327635
$ws.Range("S10").Value = -1005.576138894399
$ws.Range("T10").Value = -1005.576138894399
$ws.Range("U10").Value = -1005.576138894399
$ws.Range("V10").Value = -984.7614892657341
$ws.Range("W10").Value = -984.7614892657341
$ws.Range("X10").Value = -984.7614892657341
$ws.Range("Y10").Value = -984.7614892657341
$ws.Range("Z10").Value = -984.7614892657341
$ws.Range("AA10").Value = -984
$ws.Range("AB10").Value = 'passed'
$ws.Range("AC10").Value = 'limit(add(ang_vel(x3, x2, x3, y1), sub(ang_vel(conditional(y3, x3), limit(x2, y3, y2), sub(y3, x3), ang_vel(y2, x2, x1, x1)), y1)), sin(conditional(sin(sin(y2)), y2)), conditional(protectedDiv(limit(x2, y3, y3), x3), sub(y2, y1)))'
$ws.Range("AD10").Value = 'max, tan, atan, acos, cos, asin, '
$ws.Range("AE10").Value = 'ang_vel, y2, limit, sub, y1, conditional, x2, protectedDiv, x1, y3, x3, add, sin, '
$ws.Range("AF10").Value = 'replace'

# Row 11
$ws.Range("A11").Value = -1869.666379095403
$ws.Range("B11").Value = -1858.017431513293
$ws.Range("C11").Value = -1683.847603090057
$ws.Range("D11").Value = -1638.039076445267
$ws.Range("E11").Value = -1638.039076445267
$ws.Range("F11").Value = -1638.039076445267
$ws.Range("G11").Value = -1561.561743484393
$ws.Range("H11").Value = -1503.110190170663
$ws.Range("I11").Value = -1434.012662221568
$ws.Range("J11").Value = -1441.600235894704
$ws.Range("K11").Value = -1416.57105079192
$ws.Range("L11").Value = -1403.018450573949
$ws.Range("M11").Value = -1363.708964168934
$ws.Range("N11").Value = -1265.797229648083
$ws.Range("O11").Value = -1265.797229648083
$ws.Range("P11").Value = -1265.797229648083
$ws.Range("Q11").Value = -1224.686774815917
$ws.Range("R11").Value = -1224.686774815917
$ws.Range("S11").Value = -1224.686774815917
$ws.Range("T11").Value = -1224.686774815917
$ws.Range("U11").Value = -1224.686774815917
$ws.Range("V11").Value = -1224.686774815917
$ws.Range("W11").Value = -1224.686774815917
$ws.Range("X11").Value = -1115.965253686046
$ws.Range("Y11").Value = -1196.46924446723
$ws.Range("Z11").Value = -1221.312058234421
$ws.Range("AA11").Value = -1115
$ws.Range("AB11").Value = 'passed'
$ws.Range("AC11").Value = 'limit(atan(acos(x1, asin(x3, y1)), x2), asin(x2, x3), limit(x1, x2, limit(x1, x2, tan(sub(sub(tan(x1), max(y2, y3)), x3)))))'
$ws.Range("AD11").Value = 'ang_vel, protectedDiv, conditional, add, cos, sin, '
$ws.Range("AE11").Value = 'max, tan, y2, limit, atan, y1, sub, acos, x2, x1, y3, x3, asin, '
$ws.Range("AF11").Value = 'replace'

# Row 12
$ws.Range("A12").Value = -1734.229883634241
$ws.Range("B12").Value = -1606.824069959779
$ws.Range("C12").Value = -863.7435572017827
$ws.Range("D12").Value = -844.8124702330682
$ws.Range("E12").Value = -747.8722419058496
$ws.Range("F12").Value = -747.8722419058496
$ws.Range("G12").Value = -747.8722419058496
$ws.Range("H12").Value = -747.8722419058496
$ws.Range("I12").Value = -657.8122769901919
$ws.Range("J12").Value = -657.8122769901919
$ws.Range("K12").Value = -633.5609828241129
$ws.Range("L12").Value = -633.5609828241129
$ws.Range("M12").Value = -633.214810131536
$ws.Range("N12").Value = -548.7929538670101
$ws.Range("O12").Value = -548.7929538670101
$ws.Range("P12").Value = -548.7929538670101
$ws.Range("Q12").Value = -513.0088329597651
$ws.Range("R12").Value = -453.6275175247511
$ws.Range("S12").Value = -453.6275175247511
$ws.Range("T12").Value = -453.6275175247511
$ws.Range("U12").Value = -490.27208051843
$ws.Range("V12").Value = -490.27208051843
$ws.Range("W12").Value = -490.27208051843
$ws.Range("X12").Value = -490.27208051843
$ws.Range("Y12").Value = -490.27208051843
$ws.Range("Z12").Value = -465.3125497698164
$ws.Range("AA12").Value = -453
$ws.Range("AB12").Value = 'passed'
$ws.Range("AC12").Value = 'sub(x3, ang_vel(tan(x3), asin(x1, limit(add(acos(add(y3, x1), add(add(y2, y3), tan(y1))), y3), acos(y3, x3), asin(y1, x2))), max(tan(protectedDiv(x1, y1)), acos(x1, y2)), x1))'
$ws.Range("AD12").Value = 'atan, conditional, cos, sin, '
$ws.Range("AE12").Value = 'max, tan, ang_vel, y2, sub, limit, y1, protectedDiv, acos, x2, x1, y3, x3, add, asin, '
$ws.Range("AF12").Value = 'replace'

Write-Output "Added rows 8-12 with tree pdf data for best individual"
